$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, matching the style of the existing header row (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 10:50:59.984323"
$ws.Range("F3").Value = "2021-10-05 10:50:59.984336"
$ws.Range("F4").Value = "2021-10-05 10:50:59.984340"
$ws.Range("F5").Value = "2021-10-05 10:50:59.984344"
$ws.Range("F6").Value = "2021-10-05 10:50:59.984347"
$ws.Range("F7").Value = "2021-10-05 10:50:59.984350"
